# Generate Report for Handoff
# Updates the localization-status workbook: the first pair of files
# (b831ba9f-...md / f365c08c-...md) finished handback and a brand new
# pair of files (8f446965-...md / ffff45634dd0-...md) is now queued up for
# handoff, so the report is regenerated to reflect the new, single
# in-flight handoff (only one "generation" of target files per sheet instead
# of two).

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

$zhXlf = "8f446965-5bfd-4271-9b70-f4003cd06353.c1c80d228948d19e18cef5efebf0714bde1b1704.zh-cn.xlf"
$deXlf = "8f446965-5bfd-4271-9b70-f4003cd06353.c1c80d228948d19e18cef5efebf0714bde1b1704.de-de.xlf"
$md1 = "8f446965-5bfd-4271-9b70-f4003cd06353.md"
$md2 = "ffff45634dd0-d355-4a51-a573-a5948a724b39.md"

$md1Url = "https://github.com/OpenLocalizationTest/oltest/blob/a19b38252c42142d4bef91d170322f4458b09b54/e2e/$md1"
$md2Url = "https://github.com/OpenLocalizationTest/oltest/blob/a19b38252c42142d4bef91d170322f4458b09b54/e2e/$md2"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b781e7a3b25ab459a410f23718678d4356e89a51/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1818f284c0db9841dcd99bb91531fc8e9fbbd47e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"

# ===================== Sheet 1: Overview =====================
$ws1 = $wb.Worksheets.Item("Overview")

# existing hyperlinks need to be rebuilt (loaded hyperlinks can't be edited
# in place), so drop them all and re-add the ones we still need afterwards
$ws1.Range("A1:D3").Hyperlinks.Delete()

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-03-21 19:03:51"

$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-21 19:03:51"

$ws1.Hyperlinks.Add($ws1.Range("A2"), $md1Url, $missing, $missing, $md1)
$ws1.Hyperlinks.Add($ws1.Range("A3"), $md2Url, $missing, $missing, $md2)

# ===================== Sheet 2: zh-cn =====================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A1:L3").Hyperlinks.Delete()

# the "Latest Target File" (F) / "Latest Handback File" (G) columns no
# longer apply to a freshly-queued handoff, so clear them out entirely
$ws2.Range("F2:G3").Clear()

$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("E2").Value = "2016-03-21 19:03:46"
$ws2.Range("H2").Value = "0001-01-01 00:00:00"

$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("E3").Value = "2016-03-21 19:03:46"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"

$ws2.Hyperlinks.Add($ws2.Range("A2"), $md1Url, $missing, $missing, $md1)
$ws2.Hyperlinks.Add($ws2.Range("D2"), $zhXlfUrl, $missing, $missing, $zhXlf)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $md2Url, $missing, $missing, $md2)
$ws2.Hyperlinks.Add($ws2.Range("D3"), $zhXlfUrl, $missing, $missing, $zhXlf)

# ===================== Sheet 3: de-de =====================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A1:L3").Hyperlinks.Delete()

$ws3.Range("F2:G3").Clear()

$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("E2").Value = "2016-03-21 19:03:51"
$ws3.Range("H2").Value = "0001-01-01 00:00:00"

$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("E3").Value = "2016-03-21 19:03:51"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"

$ws3.Hyperlinks.Add($ws3.Range("A2"), $md1Url, $missing, $missing, $md1)
$ws3.Hyperlinks.Add($ws3.Range("D2"), $deXlfUrl, $missing, $missing, $deXlf)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $md2Url, $missing, $missing, $md2)
$ws3.Hyperlinks.Add($ws3.Range("D3"), $deXlfUrl, $missing, $missing, $deXlf)

Write-Host "Handoff report regenerated"
